$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Price-scrape update: append the newly scraped row (row 38) with
# Date / Price / Discount / Incredible values for 2026-02-07.
#
# The sheet stores every value (even numeric-looking ones) as a shared
# string, not as a number. Assigning plain strings via .Value would let
# Excel "smart convert" them into real numbers/dates (and pick up extra
# number-format styles), so instead we build the values with literal
# string formulas (="...") in a scratch row, copy them, and paste back
# as values only. That keeps the cells as plain text/shared-strings
# without touching the style table, exactly like the rest of the sheet.

$scratchRow = 100
$targetRow = 38

$ws.Range("A$scratchRow").Formula = '="2026-02-07"'
$ws.Range("B$scratchRow").Formula = '="24024100"'
$ws.Range("C$scratchRow").Formula = '="0"'
$ws.Range("D$scratchRow").Formula = '="0"'

$ws.Range("A$scratchRow`:D$scratchRow").Copy()
$ws.Range("A$targetRow").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Rows.Item($scratchRow).Delete()
